$wb = $excel.ActiveWorkbook

$wsGroups   = $wb.Worksheets.Item("Groups")
$wsHosts    = $wb.Worksheets.Item("Hosts")
$wsNetworks = $wb.Worksheets.Item("Networks")

# ---------------------------------------------------------------------------
# Hosts sheet: split the single "IP-Address" column into separate IPv4 /
# IPv6 columns (two new columns inserted right after the existing
# IP-Address column).
# ---------------------------------------------------------------------------
$wsHosts.Columns("C:D").Insert()

$wsHosts.Range("C1").Value = "IPv4-Address"
$wsHosts.Range("D1").Value = "IPv6-Address"

# Row 3 (TestHostB): its IP address was stored under the old "B" column;
# move it to the new IPv4-Address column.
$wsHosts.Range("C3").Value = $wsHosts.Range("B3").Value()
$wsHosts.Range("B3").ClearContents()

# Row 4 (TestHostC): same move, plus a brand new IPv6 address.
$wsHosts.Range("C4").Value = $wsHosts.Range("B4").Value()
$wsHosts.Range("B4").ClearContents()
$wsHosts.Range("D4").Value = "fd00:1234::1"

# Row 5 (TestHostD): its old IPv4 address is replaced by an IPv6 address.
$wsHosts.Range("B5").Value = "fd00:1234::2"

$wsHosts.Columns("B:D").ColumnWidth = 14.0

# ---------------------------------------------------------------------------
# Networks sheet: add SubnetMask plus separate Subnet4/MaskLength4 and
# Subnet6/MaskLength6 columns (five new columns inserted after MaskLength).
# ---------------------------------------------------------------------------
$wsNetworks.Columns("D:H").Insert()

$wsNetworks.Range("D1").Value = "SubnetMask"

# Row 3 (TestNetworkB): replace the plain subnet mask length with an
# explicit dotted-decimal SubnetMask plus an IPv6 subnet/length pair.
$wsNetworks.Range("C3").ClearContents()
$wsNetworks.Range("D3").Value = "255.255.255.128"

$wsNetworks.Range("E1").Value = "Subnet4"
$wsNetworks.Range("F1").Value = "MaskLength4"
$wsNetworks.Range("G1").Value = "Subnet6"
$wsNetworks.Range("H1").Value = "MaskLength6"

$wsNetworks.Range("G3").Value = "fd00::"
$wsNetworks.Range("H3").Value = 8

# Row 4 (TestNetworkC): its old Subnet/MaskLength move into the new
# Subnet4/MaskLength4 columns.
$wsNetworks.Range("E4").Value = $wsNetworks.Range("B4").Value()
$wsNetworks.Range("B4").ClearContents()
$wsNetworks.Range("F4").Value = $wsNetworks.Range("C4").Value()
$wsNetworks.Range("C4").ClearContents()

$wsNetworks.Columns("D:H").ColumnWidth = 14.0

# ---------------------------------------------------------------------------
# View state: Groups was the active tab before, Networks is now the active
# tab, with a new selection on Hosts and Networks as well.
# ---------------------------------------------------------------------------
$wsGroups.Range("A3").Select()
$wsHosts.Range("B5").Select()
$wsNetworks.Range("G3").Select()
$wsNetworks.Activate()
